$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 31 (shifts existing rows 31..64 down to 32..65)
# to add the "Journal of Information Technology & Politics" (JITP) entry.
$ws.Rows.Item(31).Insert()

$ws.Cells.Item(31, 1).Value = "Journal of Information Technology & Politics"
$ws.Cells.Item(31, 2).Value = "<a href='https://www.tandfonline.com/action/authorSubmission?show=instructions&journalCode=witp20'target='_blank'>Research Note</a>"
$ws.Cells.Item(31, 3).Value = "5k words"
$ws.Cells.Item(31, 4).Value = 29

# Restore the selection to the first blank row below the (now longer) table.
$ws.Range("A66:XFD66").Select()
